# Form-3-Schedule-of-Outputs.docx — "fix generate evaluation form"
#
# 1) Every table in the document should use a fixed column layout
#    (adds <w:tblLayout w:type="fixed"/> to each table's tblPr) instead
#    of Word's default auto-fit-to-contents behaviour, so the generated
#    form keeps its authored column widths.
# 2) The {{ProjectTitle}} merge field in the header table was split by
#    Word's spell-checker into three runs — "{{", "ProjectTitle", "}}" —
#    with <w:proofErr/> spell-check markers in between. A merge engine
#    that only scans whole runs for "{{Field}}" tokens would miss this
#    one, so collapse it back into a single run/token.

$d = $word.ActiveDocument

# --- 1) Fixed table layout on every table -----------------------------
$tableCount = $d.Tables.Count
for ($i = 1; $i -le $tableCount; $i++) {
    $table = $d.Tables.Item($i)
    $table.AllowAutoFit = $false
}

# --- 2) Re-merge the split {{ProjectTitle}} placeholder ----------------
$d.Content.Find.Execute(
    "{{ProjectTitle}}",  # FindText
    $true,               # MatchCase
    $false,              # MatchWholeWord
    $false,              # MatchWildcards
    $false,              # MatchSoundsLike
    $false,              # MatchAllWordForms
    $true,                # Forward
    1,                    # Wrap (wdFindContinue)
    $false,               # Format
    "{{ProjectTitle}}",   # ReplaceWith
    2                     # Replace (wdReplaceAll)
) | Out-Null
